$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.274.86'
$ws.Range('E2').Value = '  +2.72%  '
$ws.Range('D3').Value = '2.597.72'
$ws.Range('E3').Value = '  +1.61%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.83%  '
$ws.Range('D9').Value = '2.612.45'
$ws.Range('E9').Value = '  +2.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('E11').Value = '  +4.26%  '
$ws.Range('E12').Value = '  +3.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.136'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.74%  '
$ws.Range('D14').Value = '3.058.03'
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').Value = '59.219.09'
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.57'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.61%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.575.88'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '346.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.62%  '
$ws.Range('E20').Value = '  +2.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.09'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.67%  '
$ws.Range('E25').Value = '  +3.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.407'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.05%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.05%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0737'
$ws.Range('E30').Value = '  +5.53%  '
$ws.Range('E31').Value = '  +5.62%  '
$ws.Range('E32').Value = '  -1.00%  '
$ws.Range('E33').Value = '  +2.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '149.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('E35').Value = '  +3.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.08%  '
$ws.Range('E38').Value = '  +5.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.841'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.831'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('E41').Value = '  +2.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '278.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('E44').Value = '  +3.70%  '
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0962'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.98%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.943.85'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0223'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.79%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.87%  '
